$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 26, pushing the existing rows 26-113
# down to 28-115 (same formatting as the row above is copied automatically).
$ws.Range("A26:A27").EntireRow.Insert()

# --- Fill in the new row 26 ---
$ws.Range("A26").Value = 1
$ws.Range("B26").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C26").Value = "Arica y Parinacota"
$ws.Range("D26").Value = 44910
$ws.Range("E26").Value = 15
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100102
$ws.Range("H26").Value = "Cítricos"
$ws.Range("I26").Value = 100102005
$ws.Range("J26").Value = "Naranja"
$ws.Range("K26").Value = "Midknight"
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 200
$ws.Range("N26").Value = 950
$ws.Range("O26").Value = 1000
$ws.Range("P26").Value = 975
$ws.Range("Q26").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R26").Value = "Región de O'Higgins"
$ws.Range("S26").Value = 975
$ws.Range("T26").Value = 1

# --- Fill in the new row 27 ---
$ws.Range("A27").Value = 1
$ws.Range("B27").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C27").Value = "Arica y Parinacota"
$ws.Range("D27").Value = 44910
$ws.Range("E27").Value = 15
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100102
$ws.Range("H27").Value = "Cítricos"
$ws.Range("I27").Value = 100102005
$ws.Range("J27").Value = "Naranja"
$ws.Range("K27").Value = "Midknight"
$ws.Range("L27").Value = "Tercera"
$ws.Range("M27").Value = 180
$ws.Range("N27").Value = 800
$ws.Range("O27").Value = 850
$ws.Range("P27").Value = 833
$ws.Range("Q27").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R27").Value = "Región de O'Higgins"
$ws.Range("S27").Value = 833
$ws.Range("T27").Value = 1
